$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear B2 and C2 (previously "devserver" and " ")
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = $null

# A2 keeps its text/hyperlink as-is (http://172.191.4.85/TestCollection)

# Row 3: apply style changes -- A3 gets the Hyperlink style (like A2), B3 gets wrap-text style (like B2)
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("B3").WrapText = $true

# Update selection to B2
$ws.Range("B2").Select()

